$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-29 hold the per-trial stimulus list for this block. The list was
# regenerated (re-randomized old/new/target assignment + renumbered trial_total)
# and the catch-trial image/position was swapped. Rewrite each row A:S in place.

$row2 = @(5, "memory", 6, 2, 1, 297, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_vh7v8.png", 78.70454545454545, 59.63636363636363, 69.17045454545455, 44, 7, 7, 7)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }

$row3 = @(5, "memory", 6, 2, 2, 298, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_xr3up.png", 76.24444444444444, 55.88888888888889, 66.06666666666666, 45, 7, 7, 7)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }

$row4 = @(5, "memory", 6, 2, 3, 299, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_ra2nm.png", 70.75, 50.375, 60.5625, 40, 6, 6, 6)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }

$row5 = @(5, "memory", 6, 2, 4, 300, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_9oofc.png", 82.47619047619048, 65.5, 73.98809523809524, 42, 8, 8, 8)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }

$row6 = @(5, "memory", 6, 2, 5, 301, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_syam3.png", 41.32432432432432, 26.2972972972973, 33.81081081081081, 37, 2, 2, 2)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row6[$i] }

$row7 = @(5, "memory", 6, 2, 6, 302, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_tn8ys.png", 86.70454545454545, 72.4090909090909, 79.55681818181819, 44, 10, 10, 10)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row7[$i] }

$row8 = @(5, "memory", 6, 2, 7, 303, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_bf8nx.png", 86.63414634146342, 66.63414634146342, 76.63414634146342, 41, 9, 9, 9)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row8[$i] }

$row9 = @(5, "memory", 6, 2, 8, 304, "living_rooms", $null, $null, "catch", "f", "stimuli/catch_06.jpg", $null, $null, $null, $null, $null, $null, $null)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row9[$i] }

$row10 = @(5, "memory", 6, 2, 9, 305, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_7lz7m.png", 51.5531914893617, 32.87234042553192, 42.21276595744681, 47, 3, 3, 3)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row10[$i] }

$row11 = @(5, "memory", 6, 2, 10, 306, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_2qhro.png", 81.73809523809524, 62.73809523809524, 72.23809523809524, 42, 8, 8, 8)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row11[$i] }

$row12 = @(5, "memory", 6, 2, 11, 307, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_ac0ey.png", 86.62222222222222, 70.02222222222223, 78.32222222222222, 45, 9, 9, 9)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row12[$i] }

$row13 = @(5, "memory", 6, 2, 12, 308, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_rru0v.png", 56.45238095238095, 39.42857142857143, 47.94047619047619, 42, 4, 4, 4)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row13[$i] }

$row14 = @(5, "memory", 6, 2, 13, 309, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_of8d6.png", 26.04878048780488, 19.14634146341463, 22.59756097560975, 41, 1, 1, 1)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row14[$i] }

$row15 = @(5, "memory", 6, 2, 14, 310, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_vgh2g.png", 93.81395348837209, 78.27906976744185, 86.04651162790697, 43, 10, 10, 10)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row15[$i] }

$row16 = @(5, "memory", 6, 2, 15, 311, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_kq9s9.png", 62.30232558139535, 39.97674418604651, 51.13953488372093, 43, 4, 4, 4)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row16[$i] }

$row17 = @(5, "memory", 6, 2, 16, 312, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_iudc4.png", 73.625, 52.275, 62.95, 40, 6, 6, 6)
for ($i = 0; $i -lt $row17.Length; $i++) { $ws.Cells.Item(17, $i + 1).Value = $row17[$i] }

$row18 = @(5, "memory", 6, 2, 17, 313, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_qdln8.png", 85.51162790697674, 67.86046511627907, 76.68604651162791, 43, 9, 9, 9)
for ($i = 0; $i -lt $row18.Length; $i++) { $ws.Cells.Item(18, $i + 1).Value = $row18[$i] }

$row19 = @(5, "memory", 6, 2, 18, 314, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_pdzf1.png", 86.23913043478261, 67.17391304347827, 76.70652173913044, 46, 9, 9, 9)
for ($i = 0; $i -lt $row19.Length; $i++) { $ws.Cells.Item(19, $i + 1).Value = $row19[$i] }

$row20 = @(5, "memory", 6, 2, 19, 315, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_rg4in.png", 49.3695652173913, 30.21739130434782, 39.79347826086956, 46, 3, 3, 3)
for ($i = 0; $i -lt $row20.Length; $i++) { $ws.Cells.Item(20, $i + 1).Value = $row20[$i] }

$row21 = @(5, "memory", 6, 2, 20, 316, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_swq34.png", 64.11363636363636, 43.04545454545455, 53.57954545454545, 44, 5, 5, 5)
for ($i = 0; $i -lt $row21.Length; $i++) { $ws.Cells.Item(21, $i + 1).Value = $row21[$i] }

$row22 = @(5, "memory", 6, 2, 21, 317, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_3jnt7.png", 49.52272727272727, 35.25, 42.38636363636364, 44, 3, 3, 3)
for ($i = 0; $i -lt $row22.Length; $i++) { $ws.Cells.Item(22, $i + 1).Value = $row22[$i] }

$row23 = @(5, "memory", 6, 2, 22, 318, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_pna7l.png", 85.53333333333333, 67.97777777777777, 76.75555555555556, 45, 9, 9, 9)
for ($i = 0; $i -lt $row23.Length; $i++) { $ws.Cells.Item(23, $i + 1).Value = $row23[$i] }

$row24 = @(5, "memory", 6, 2, 23, 319, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_amsgw.png", 86.08510638297872, 65.95744680851064, 76.02127659574468, 47, 9, 9, 9)
for ($i = 0; $i -lt $row24.Length; $i++) { $ws.Cells.Item(24, $i + 1).Value = $row24[$i] }

$row25 = @(5, "memory", 6, 2, 24, 320, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_rych7.png", 30.4468085106383, 23.4468085106383, 26.9468085106383, 47, 2, 2, 2)
for ($i = 0; $i -lt $row25.Length; $i++) { $ws.Cells.Item(25, $i + 1).Value = $row25[$i] }

$row26 = @(5, "memory", 6, 2, 25, 321, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_zxvl3.png", 68.78260869565217, 47.56521739130435, 58.17391304347827, 46, 5, 5, 5)
for ($i = 0; $i -lt $row26.Length; $i++) { $ws.Cells.Item(26, $i + 1).Value = $row26[$i] }

$row27 = @(5, "memory", 6, 2, 26, 322, "living_rooms", "living_rooms", $null, "new", "f", "stimuli/img_lgxzn.png", 73.11363636363636, 49.97727272727273, 61.54545454545455, 44, 6, 6, 6)
for ($i = 0; $i -lt $row27.Length; $i++) { $ws.Cells.Item(27, $i + 1).Value = $row27[$i] }

$row28 = @(5, "memory", 6, 2, 27, 323, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_j4ttn.png", 12.61904761904762, 11.42857142857143, 12.02380952380952, 42, 1, 1, 1)
for ($i = 0; $i -lt $row28.Length; $i++) { $ws.Cells.Item(28, $i + 1).Value = $row28[$i] }

$row29 = @(5, "memory", 6, 2, 28, 324, "living_rooms", "living_rooms", "target", "old", "j", "stimuli/img_24rt2.png", 55.26829268292683, 34.19512195121951, 44.73170731707317, 41, 3, 3, 3)
for ($i = 0; $i -lt $row29.Length; $i++) { $ws.Cells.Item(29, $i + 1).Value = $row29[$i] }

